# Monopoly_data.xlsx — "Completed all cards classes and populate in dictionaries"
#
# The generic "Community" / "Chance" labels used in column A of the
# "Community Cards" / "Chance Cards" sheets are replaced with per-row
# labels ("Community1".."Community10" / "Chance1".."Chance10"), and the
# generic "City" header used atop the first column of the "Special Cards",
# "Transportation" and "Energy" sheets is replaced with a sheet-specific
# header ("Block" / "Route" / "Station").
#
# NOTE on ordering: the shared-strings table is rebuilt on save with
# unused strings pruned and new strings appended in first-use order, so
# the order in which we touch the sheets below controls the resulting
# shared-string indices. Touch them in the same order the new strings
# appear in the target workbook: Block, Route, Station, Chance1..10,
# Community1..10.

$wb = $excel.ActiveWorkbook

$wsSpecial = $wb.Worksheets.Item("Special Cards")
$wsSpecial.Range("A1").Value = "Block"

$wsTransport = $wb.Worksheets.Item("Transportation")
$wsTransport.Range("A1").Value = "Route"

$wsEnergy = $wb.Worksheets.Item("Energy")
$wsEnergy.Range("A1").Value = "Station"

$wsChance = $wb.Worksheets.Item("Chance Cards")
for ($r = 2; $r -le 11; $r++) {
    $wsChance.Cells.Item($r, 1).Value = "Chance" + ($r - 1)
}

$wsCommunity = $wb.Worksheets.Item("Community Cards")
for ($r = 2; $r -le 11; $r++) {
    $wsCommunity.Cells.Item($r, 1).Value = "Community" + ($r - 1)
}

# Widen column A on the "Chance Cards" sheet (it now holds "Chance1".."Chance10"
# labels instead of the single repeated "Chance" label).
$wsChance.Columns("A").ColumnWidth = 29

# Restore the selected cell on every touched sheet (setting .Value above does
# not move the selection, but re-assert explicitly for clarity/robustness).
# Sheets are selected in order so that the sheet that should end up active
# ("Community Cards", which was already the active tab before the edit) is
# selected last.
$wsSpecial.Range("A2").Select() | Out-Null
$wsTransport.Range("E16").Select() | Out-Null
$wsEnergy.Range("A4").Select() | Out-Null
$wsChance.Range("B15").Select() | Out-Null
$wsCommunity.Range("B19").Select() | Out-Null
